# Update the NATMI lrc2p edge-weight table with re-computed TPM values.
#
# The "Ligand average/total expression value" (G/H) for the ECs and MuSCs
# sending clusters changed (new TPM normalization); the FAPs sending
# cluster's G/H is untouched. Everything else in the row is derived:
#   I = G / sum(G over the 3 sending clusters)   (ligand avg specificity)
#   J = H / sum(H over the 3 sending clusters)   (ligand total specificity)
#   Q = G * M   (edge average weight = ligand avg * receptor avg)
#   R = H * N   (edge total weight  = ligand total * receptor total)
#   S = I * O   (edge average specificity = ligand avg spec * receptor avg spec)
#   T = J * P   (edge total specificity  = ligand total spec * receptor total spec)
#
# Because the normalizing sum (over all 3 sending clusters) changes, I/J/S/T
# change for every row -- even rows whose own G/H (FAPs) did not move. Q/R
# only depend on that row's own G/H, so they are only rewritten for rows
# whose sending cluster's G/H actually changed (ECs, MuSCs), leaving the
# FAPs rows' Q/R byte-identical to the original.
#
# NOTE: reading back a cell's `.Value` right after assignment can return the
# property descriptor instead of the number in this host, so `.Value2` is
# used for every read (it behaves correctly for both get and set).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New raw ligand average/total expression values for the sending clusters
# whose underlying TPM changed (FAPs is untouched).
$newG = @{ "ECs" = 0.8986206666666666; "FAPs" = 2.109481;          "MuSCs" = 2.399397 }
$newH = @{ "ECs" = 2.695862;           "FAPs" = 6.328443;          "MuSCs" = 7.198191 }

# Data rows: row number -> sending cluster (col A) label used as key above.
$sendingCluster = @{ 2 = "ECs"; 3 = "ECs"; 4 = "FAPs"; 5 = "FAPs"; 6 = "MuSCs"; 7 = "MuSCs" }

# Rows whose own ligand average/total expression value actually changed.
$changedGHRows = @(2, 3, 6, 7)

# 1) Write the updated G/H (ligand average / total expression value) cells
#    only where the underlying value actually changed.
foreach ($r in $changedGHRows) {
    $cluster = $sendingCluster[$r]
    $ws.Cells.Item($r, 7).Value2 = $newG[$cluster]   # column G
    $ws.Cells.Item($r, 8).Value2 = $newH[$cluster]   # column H
}

# 2) Derived specificity of average/total ligand expression (I, J) is each
#    sending cluster's value normalized across all three sending clusters;
#    the normalizing sum changed, so every row's I/J is refreshed.
$sumG = 0.0
$sumH = 0.0
foreach ($cluster in @("ECs", "FAPs", "MuSCs")) {
    $sumG += $newG[$cluster]
    $sumH += $newH[$cluster]
}

foreach ($r in 2..7) {
    $cluster = $sendingCluster[$r]
    $ws.Cells.Item($r, 9).Value2  = $newG[$cluster] / $sumG   # column I
    $ws.Cells.Item($r, 10).Value2 = $newH[$cluster] / $sumH   # column J
}

# 3) Edge average/total weight (Q, R) only change where G/H changed.
foreach ($r in $changedGHRows) {
    $g = $ws.Cells.Item($r, 7).Value2
    $h = $ws.Cells.Item($r, 8).Value2
    $m = $ws.Cells.Item($r, 13).Value2
    $n = $ws.Cells.Item($r, 14).Value2

    $ws.Cells.Item($r, 17).Value2 = $g * $m   # column Q
    $ws.Cells.Item($r, 18).Value2 = $h * $n   # column R
}

# 4) Edge average/total specificity (S, T) depend on I/J, which changed for
#    every row, so every row's S/T is refreshed.
foreach ($r in 2..7) {
    $i = $ws.Cells.Item($r, 9).Value2
    $j = $ws.Cells.Item($r, 10).Value2
    $o = $ws.Cells.Item($r, 15).Value2
    $p = $ws.Cells.Item($r, 16).Value2

    $ws.Cells.Item($r, 19).Value2 = $i * $o   # column S
    $ws.Cells.Item($r, 20).Value2 = $j * $p   # column T
}
